$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.159054695133975
$ws.Range("C2").Value = 4.500813361879064
$ws.Range("E2").Value = 20.84741186657752
$ws.Range("F2").Value = 40.43597167435549
$ws.Range("G2").Value = 3.64750983676834
$ws.Range("I2").Value = 21.82737116811017
$ws.Range("J2").Value = 8.305246227963558
$ws.Range("K2").Value = 8.131938830123946
$ws.Range("M2").Value = 17.63623549835984
$ws.Range("O2").Value = 23.1165425373109
$ws.Range("B3").Value = 7.856956119717805
$ws.Range("C3").Value = 4.327538061331373
$ws.Range("E3").Value = 20.73292576435665
$ws.Range("F3").Value = 40.40830126207543
$ws.Range("G3").Value = 3.649135742441822
$ws.Range("I3").Value = 21.93765498464185
$ws.Range("J3").Value = 8.320482328259299
$ws.Range("K3").Value = 7.935173013125079
$ws.Range("M3").Value = 17.50171694112018
$ws.Range("O3").Value = 23.22213170804183
$ws.Range("B4").Value = 7.666183235077403
$ws.Range("C4").Value = 4.216721847421771
$ws.Range("E4").Value = 20.66658662842658
$ws.Range("F4").Value = 40.40161096148864
$ws.Range("G4").Value = 3.650186352986709
$ws.Range("I4").Value = 22.00946642769632
$ws.Range("J4").Value = 8.330318766210246
$ws.Range("K4").Value = 7.812678738067859
$ws.Range("M4").Value = 17.42147761401865
$ws.Range("O4").Value = 23.29175205555968
$ws.Range("B5").Value = 7.58722835054327
$ws.Range("C5").Value = 4.170492182668474
$ws.Range("E5").Value = 20.64056879463397
$ws.Range("F5").Value = 40.40147726440889
$ws.Range("G5").Value = 3.650627676789443
$ws.Range("I5").Value = 22.03976089017069
$ws.Range("J5").Value = 8.334448605126152
$ws.Range("K5").Value = 7.762407816875236
$ws.Range("M5").Value = 17.38939982199019
$ws.Range("O5").Value = 23.32132563297427
$ws.Range("B6").Value = 7.574048266503117
$ws.Range("C6").Value = 4.162752402013241
$ws.Range("E6").Value = 20.63631050066853
$ws.Range("F6").Value = 40.40161169700116
$ws.Range("G6").Value = 3.650701756202125
$ws.Range("I6").Value = 22.04485353265091
$ws.Range("J6").Value = 8.335141705104384
$ws.Range("K6").Value = 7.754041090867612
$ws.Range("M6").Value = 17.38411159782061
$ws.Range("O6").Value = 23.32630890691263
$ws.Range("B7").Value = 7.665123173916963
$ws.Range("C7").Value = 4.216102656439022
$ws.Range("E7").Value = 20.66623160261548
$ws.Range("F7").Value = 40.4015986586629
$ws.Range("G7").Value = 3.650192251367534
$ws.Range("I7").Value = 22.00987081531745
$ws.Range("J7").Value = 8.33037397060369
$ws.Range("K7").Value = 7.812002104939538
$ws.Range("M7").Value = 17.42104245351649
$ws.Range("O7").Value = 23.29214602711178
$ws.Range("B8").Value = 8.056052987520882
$ws.Range("C8").Value = 4.442009564668836
$ws.Range("E8").Value = 20.8071291917035
$ws.Range("F8").Value = 40.42429550401318
$ws.Range("G8").Value = 3.648059619197098
$ws.Range("I8").Value = 21.86454708246615
$ws.Range("J8").Value = 8.310399933416742
$ws.Range("K8").Value = 8.064482717676006
$ws.Range("M8").Value = 17.58938172586954
$ws.Range("O8").Value = 23.15195511066727
$ws.Range("B9").Value = 8.776286318001116
$ws.Range("C9").Value = 4.848345383186716
$ws.Range("E9").Value = 21.11381854304573
$ws.Range("F9").Value = 40.55032105585345
$ws.Range("G9").Value = 3.644290629963493
$ws.Range("I9").Value = 21.61204781831677
$ws.Range("J9").Value = 8.275034334681006
$ws.Range("K9").Value = 8.543428424884562
$ws.Range("M9").Value = 17.93690840026691
$ws.Range("O9").Value = 22.91509174618701
$ws.Range("B10").Value = 9.271985256580345
$ws.Range("C10").Value = 5.122763680016297
$ws.Range("E10").Value = 21.35624362162552
$ws.Range("F10").Value = 40.69221210392202
$ws.Range("G10").Value = 3.641770771251804
$ws.Range("I10").Value = 21.4463015712924
$ws.Range("J10").Value = 8.251346857958247
$ws.Range("K10").Value = 8.881913790759652
$ws.Range("M10").Value = 18.20108825282544
$ws.Range("O10").Value = 22.76434571668573
$ws.Range("B11").Value = 9.489339330961348
$ws.Range("C11").Value = 5.24207625535829
$ws.Range("E11").Value = 21.46991400862835
$ws.Range("F11").Value = 40.76734100666726
$ws.Range("G11").Value = 3.640677980617394
$ws.Range("I11").Value = 21.3751845310387
$ws.Range("J11").Value = 8.241064377419118
$ws.Range("K11").Value = 9.032322666494141
$ws.Range("M11").Value = 18.32280500204238
$ws.Range("O11").Value = 22.70084012917906
$ws.Range("B12").Value = 9.57041329658793
$ws.Range("C12").Value = 5.28644311155795
$ws.Range("E12").Value = 21.51341633450878
$ws.Range("F12").Value = 40.79729801034138
$ws.Range("G12").Value = 3.640271821262995
$ws.Range("I12").Value = 21.34886981255156
$ws.Range("J12").Value = 8.237241204805628
$ws.Range("K12").Value = 9.08871516483655
$ws.Range("M12").Value = 18.36908427005715
$ws.Range("O12").Value = 22.67752289590976
$ws.Range("B13").Value = 9.553008326792488
$ws.Range("C13").Value = 5.276924417881569
$ws.Range("E13").Value = 21.50402744941619
$ws.Range("F13").Value = 40.7907794450612
$ws.Range("G13").Value = 3.640358955003264
$ws.Range("I13").Value = 21.35450977039905
$ws.Range("J13").Value = 8.238061459507984
$ws.Range("K13").Value = 9.076595882746114
$ws.Range("M13").Value = 18.35910938769278
$ws.Range("O13").Value = 22.68251213825848
$ws.Range("B14").Value = 9.49603438085175
$ws.Range("C14").Value = 5.24574276662112
$ws.Range("E14").Value = 21.47348396854278
$ws.Range("F14").Value = 40.76977546782279
$ws.Range("G14").Value = 3.640644412403509
$ws.Range("I14").Value = 21.37300726203165
$ws.Range("J14").Value = 8.240748430030015
$ws.Range("K14").Value = 9.036973692151582
$ws.Range("M14").Value = 18.32660885528109
$ws.Range("O14").Value = 22.69890714557957
$ws.Range("B15").Value = 9.460973867530214
$ws.Range("C15").Value = 5.22653652068454
$ws.Range("E15").Value = 21.45483393191378
$ws.Range("F15").Value = 40.75710576240879
$ws.Range("G15").Value = 3.64082025933524
$ws.Range("I15").Value = 21.38441770259357
$ws.Range("J15").Value = 8.242403459053671
$ws.Range("K15").Value = 9.012629061538187
$ws.Range("M15").Value = 18.30672479313101
$ws.Range("O15").Value = 22.70904481166516
$ws.Range("B16").Value = 9.257611285671107
$ws.Range("C16").Value = 5.114853413444089
$ws.Range("E16").Value = 21.34888069334812
$ws.Range("F16").Value = 40.68751400845315
$ws.Range("G16").Value = 3.641843261165447
$ws.Range("I16").Value = 21.45103539985264
$ws.Range("J16").Value = 8.25202873638859
$ws.Range("K16").Value = 8.872007950502566
$ws.Range("M16").Value = 18.19316189087369
$ws.Range("O16").Value = 22.76859813499946
$ws.Range("B17").Value = 9.130722427074572
$ws.Range("C17").Value = 5.044909830245811
$ws.Range("E17").Value = 21.28472923787408
$ws.Range("F17").Value = 40.64752244268565
$ws.Range("G17").Value = 3.642484517108405
$ws.Range("I17").Value = 21.49299987584867
$ws.Range("J17").Value = 8.258059588685065
$ws.Range("K17").Value = 8.78478971145957
$ws.Range("M17").Value = 18.12386525774565
$ws.Range("O17").Value = 22.80643198324791
$ws.Range("B18").Value = 9.056976134690647
$ws.Range("C18").Value = 5.004161402413667
$ws.Range("E18").Value = 21.24815193279999
$ws.Range("F18").Value = 40.62551736347951
$ws.Range("G18").Value = 3.642858388990446
$ws.Range("I18").Value = 21.51753974498049
$ws.Range("J18").Value = 8.26157480479219
$ws.Range("K18").Value = 8.734290846970971
$ws.Range("M18").Value = 18.0841546503639
$ws.Range("O18").Value = 22.82867006114709
$ws.Range("B19").Value = 9.031877874128826
$ws.Range("C19").Value = 4.990276259749292
$ws.Range("E19").Value = 21.23582349342476
$ws.Range("F19").Value = 40.61823846380459
$ws.Range("G19").Value = 3.642985842081824
$ws.Range("I19").Value = 21.52591774033923
$ws.Range("J19").Value = 8.262772980759541
$ws.Range("K19").Value = 8.717137179203469
$ws.Range("M19").Value = 18.07073558888497
$ws.Range("O19").Value = 22.83628136676019
$ws.Range("B20").Value = 9.144309428020211
$ws.Range("C20").Value = 5.052409286748408
$ws.Range("E20").Value = 21.2915252685358
$ws.Range("F20").Value = 40.6516765237091
$ws.Range("G20").Value = 3.642415733131602
$ws.Range("I20").Value = 21.48849097528702
$ws.Range("J20").Value = 8.257412791246294
$ws.Range("K20").Value = 8.794109111771578
$ws.Range("M20").Value = 18.13122703388341
$ws.Range("O20").Value = 22.80235511862339
$ws.Range("B21").Value = 9.512802938584979
$ws.Range("C21").Value = 5.254923817361664
$ws.Range("E21").Value = 21.48244314455853
$ws.Range("F21").Value = 40.77590406121851
$ws.Range("G21").Value = 3.640560359139223
$ws.Range("I21").Value = 21.3675573880622
$ws.Range("J21").Value = 8.239957288708876
$ws.Range("K21").Value = 9.048627382743627
$ws.Range("M21").Value = 18.33615022226151
$ws.Range("O21").Value = 22.69407168079405
$ws.Range("B22").Value = 9.746422778692509
$ws.Range("C22").Value = 5.382525531164517
$ws.Range("E22").Value = 21.60987323401766
$ws.Range("F22").Value = 40.86587168103038
$ws.Range("G22").Value = 3.639392377232653
$ws.Range("I22").Value = 21.29210949355627
$ws.Range("J22").Value = 8.228960352515537
$ws.Range("K22").Value = 9.211661275037132
$ws.Range("M22").Value = 18.47115820244275
$ws.Range("O22").Value = 22.62756369345639
$ws.Range("B23").Value = 9.622413920097781
$ws.Range("C23").Value = 5.314862925020472
$ws.Range("E23").Value = 21.54162853318207
$ws.Range("F23").Value = 40.81705636313168
$ws.Range("G23").Value = 3.64001168138543
$ws.Range("I23").Value = 21.33204901035592
$ws.Range("J23").Value = 8.234792101665587
$ws.Range("K23").Value = 9.124965443566335
$ws.Range("M23").Value = 18.39901433954531
$ws.Range("O23").Value = 22.66266967274092
$ws.Range("B24").Value = 9.138169220041195
$ws.Range("C24").Value = 5.049020453344744
$ws.Range("E24").Value = 21.28845183319365
$ws.Range("F24").Value = 40.64979538937519
$ws.Range("G24").Value = 3.642446814150448
$ws.Range("I24").Value = 21.49052815983454
$ws.Range("J24").Value = 8.257705058827916
$ws.Range("K24").Value = 8.789896915653449
$ws.Range("M24").Value = 18.12789837064957
$ws.Range("O24").Value = 22.80419675259439
$ws.Range("B25").Value = 8.58698593056109
$ws.Range("C25").Value = 4.742544835478159
$ws.Range("E25").Value = 21.02773978802499
$ws.Range("F25").Value = 40.50753354096678
$ws.Range("G25").Value = 3.645266288246106
$ws.Range("I25").Value = 21.6768822606069
$ws.Range("J25").Value = 8.415972823146314
$ws.Range("K25").Value = 8.543428424884562
$ws.Range("M25").Value = 17.84120840532555
$ws.Range("O25").Value = 22.97508779287008
